# Applies the notes.docx edit:
#  1. Moves the (hidden) "_GoBack" bookmark from its old spot (a lone
#     paragraph right after "...the trace each square are identical")
#     into the middle of the word "reality" ("In rea|lity") - this is
#     simply where the author's cursor was when the file was last saved.
#  2. Fixes the grammar "are identical" -> "is identical" in the last
#     paragraph and clears the now-stale gramStart/gramEnd proofing
#     marks that used to flag "are".
#
# Only $word / $d is used, per the harness contract.

$d = $word.ActiveDocument

# --- Step 1: relocate the _GoBack bookmark -------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$findReality = $d.Content
$hasReality = $findReality.Find.Execute("In reality", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hasReality) {
    # Split "In reality" after its 6th character ("In rea" | "lity") - that
    # is where the bookmark sits in the target document.
    $splitAt = $findReality.Start + 6
    $bmRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- Step 2: "are identical" -> "is identical" ---------------------------

# Replace across the whole phrase (not just the single word) so that the
# stale gramStart/gramEnd proofErr markers bracketing "are" get swept away
# along with it.
$findPhrase = $d.Content
$hasPhrase = $findPhrase.Find.Execute("square are identical", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hasPhrase) {
    $findPhrase.Text = "square is identical"
}

# Re-split "is" back out into its own run (matching the original run
# layout, just with "are" swapped for "is") by dropping a throwaway
# bookmark at each edge - adding a bookmark forces a run boundary, and the
# split survives the bookmark's removal.
$findIs = $d.Content
$hasIs = $findIs.Find.Execute("is identical", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hasIs) {
    $isStart = $findIs.Start
    $isEnd = $findIs.Start + 2

    $edge1 = $d.Range($isStart, $isStart)
    $d.Bookmarks.Add("zzTmpSplit1", $edge1)
    $edge2 = $d.Range($isEnd, $isEnd)
    $d.Bookmarks.Add("zzTmpSplit2", $edge2)

    $d.Bookmarks("zzTmpSplit1").Delete()
    $d.Bookmarks("zzTmpSplit2").Delete()
}
